$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 297, shifting existing rows 297:312 down to 298:313
$ws.Rows.Item(297).Insert()

# Populate the newly inserted row 297 with the new record
$ws.Cells.Item(297, 1).Value = 4
$ws.Cells.Item(297, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(297, 3).Value = "Los Lagos"
$ws.Cells.Item(297, 4).Value = 44826
$ws.Cells.Item(297, 4).NumberFormat = $ws.Cells.Item(298, 4).NumberFormat
$ws.Cells.Item(297, 5).Value = 10
$ws.Cells.Item(297, 6).Value = 100112021
$ws.Cells.Item(297, 7).Value = "Ají"
$ws.Cells.Item(297, 8).Value = "Inferno"
$ws.Cells.Item(297, 9).Value = "Primera"
$ws.Cells.Item(297, 10).Value = 90
$ws.Cells.Item(297, 11).Value = 23000
$ws.Cells.Item(297, 12).Value = 23000
$ws.Cells.Item(297, 13).Value = 23000
$ws.Cells.Item(297, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(297, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(297, 16).Value = 2300
$ws.Cells.Item(297, 17).Value = 10
$ws.Cells.Item(297, 18).Value = "Hortaliza"
